# Weekly update: insert two new rows of data (new week) above the
# existing historical rows, pushing rows 4-7 down to rows 6-9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 4 (existing rows 4-7 shift down to 6-9,
# carrying their formatting/styles with them, matching Excel's default
# "insert copies formatting from the row above" behaviour).
$ws.Rows("4:5").Insert()

# New row 4: Calidad "Especial" for the new price-reporting week.
$ws.Range("A4").Value = 7
$ws.Range("B4").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C4").Value = "Ñuble"
$ws.Range("D4").Value = 45040
$ws.Range("E4").Value = 16
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100104
$ws.Range("H4").Value = "Frutos de pepita"
$ws.Range("I4").Value = 100104003
$ws.Range("J4").Value = "Membrillo"
$ws.Range("K4").Value = "Champion"
$ws.Range("L4").Value = "Especial"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13000
$ws.Range("P4").Value = 13000
$ws.Range("Q4").Value = "$/caja 18 kilos empedrada"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 722
$ws.Range("T4").Value = 18

# New row 5: Calidad "Primera" for the new price-reporting week.
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 45040
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = "Fruta"
$ws.Range("G5").Value = 100104
$ws.Range("H5").Value = "Frutos de pepita"
$ws.Range("I5").Value = 100104003
$ws.Range("J5").Value = "Membrillo"
$ws.Range("K5").Value = "Champion"
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 40
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("Q5").Value = "$/caja 18 kilos empedrada"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 667
$ws.Range("T5").Value = 18
